$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.251.36'
$ws.Range("E2").Value = '  -1.44%  '
$ws.Range("D3").Value = '2.262.69'
$ws.Range("E3").Value = '  -1.35%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").Value = '''113.81'
$ws.Range("E5").Value = '  +5.55%  '
$ws.Range("D6").Value = '''264.94'
$ws.Range("E6").Value = '  -2.38%  '
$ws.Range("D7").Value = '''0.618'
$ws.Range("E7").Value = '  -0.17%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("D9").Value = '''0.602'
$ws.Range("E9").Value = '  -1.78%  '
$ws.Range("D10").Value = '''47.73'
$ws.Range("E10").Value = '  +1.67%  '
$ws.Range("D11").Value = '''0.0925'
$ws.Range("E11").Value = '  -1.18%  '
$ws.Range("D12").Value = '''8.75'
$ws.Range("E12").Value = '  +5.11%  '
$ws.Range("E13").Value = '  -0.51%  '
$ws.Range("D14").Value = '''15.47'
$ws.Range("E14").Value = '  -0.96%  '
$ws.Range("D15").Value = '2.602.16'
$ws.Range("E15").Value = '  -1.36%  '
$ws.Range("D16").Value = '''0.853'
$ws.Range("E16").Value = '  -0.36%  '
$ws.Range("D17").Value = '2.262.22'
$ws.Range("E17").Value = '  -1.32%  '
$ws.Range("D18").Value = '43.142.32'
$ws.Range("E18").Value = '  -1.59%  '
$ws.Range("E19").Value = '  -2.57%  '
$ws.Range("D20").Value = '''6.92'
$ws.Range("E20").Value = '  +10.02%  '
$ws.Range("D21").Value = '''71.06'
$ws.Range("E21").Value = '  -1.50%  '
$ws.Range("D22").Value = '''2.40'
$ws.Range("E22").Value = '  -3.82%  '
$ws.Range("D23").Value = '''9.85'
$ws.Range("E23").Value = '  +5.80%  '
$ws.Range("D24").Value = '''230.37'
$ws.Range("E24").Value = '  -1.37%  '
$ws.Range("D25").Value = '''2.84'
$ws.Range("E25").Value = '  -3.44%  '
$ws.Range("D26").Value = '''0.999'
$ws.Range("E26").Value = '  -0.01%  '
$ws.Range("D27").Value = '''11.34'
$ws.Range("E27").Value = '  +0.03%  '
$ws.Range("D28").Value = '''3.91'
$ws.Range("E28").Value = '  -0.24%  '
$ws.Range("D29").Value = '''41.34'
$ws.Range("E29").Value = '  +1.53%  '
$ws.Range("D30").Value = '''3.38'
$ws.Range("E30").Value = '  -1.93%  '
$ws.Range("E31").Value = '  -1.65%  '
$ws.Range("D32").Value = '''171.62'
$ws.Range("E32").Value = '  -3.52%  '
$ws.Range("D33").Value = '''21.27'
$ws.Range("E33").Value = '  -2.71%  '
$ws.Range("D34").Value = '''0.0902'
$ws.Range("E34").Value = '  -0.67%  '
$ws.Range("D35").Value = '''5.56'
$ws.Range("E35").Value = '  +0.08%  '
$ws.Range("E36").Value = '  -0.74%  '
$ws.Range("D37").Value = '''4.59'
$ws.Range("E37").Value = '  -6.19%  '
$ws.Range("D38").Value = '''0.0348'
$ws.Range("E38").Value = '  -3.67%  '
$ws.Range("D39").Value = '''3.76'
$ws.Range("E39").Value = '  +3.56%  '
$ws.Range("E40").Value = '  -8.48%  '
$ws.Range("D41").Value = '''14.29'
$ws.Range("E41").Value = '  +16.54%  '
$ws.Range("D42").Value = '''75.20'
$ws.Range("E42").Value = '  +13.93%  '
$ws.Range("D43").Value = '''2.40'
$ws.Range("E43").Value = '  +2.77%  '
$ws.Range("D44").Value = '''0.234'
$ws.Range("E44").Value = '  -1.33%  '
$ws.Range("D45").Value = '''6.15'
$ws.Range("E45").Value = '  +12.27%  '
$ws.Range("E46").Value = '  -0.04%  '
$ws.Range("D47").Value = '''1.37'
$ws.Range("E47").Value = '  -0.50%  '
$ws.Range("D48").Value = '''8.56'
$ws.Range("E48").Value = '  -2.50%  '
$ws.Range("D49").Value = '''0.0987'
$ws.Range("E49").Value = '  -2.75%  '
$ws.Range("E50").Value = '  +0.02%  '
$ws.Range("D51").Value = '''99.97'
$ws.Range("E51").Value = '  +0.43%  '
